# Update Recommandations sheet (sheet 1)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the SUCRIVOIRE (SCRC) row - it drops out of the list entirely.
# Deleting row 34 shifts BERNABE CI (BNBC) up from row 35 to row 34 automatically.
$ws1.Rows.Item(34).Delete() | Out-Null

# Apply the refreshed values (sector aggregates + per-stock reshuffle from the latest BRVM data pull).

# Row 2
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 721.65

# Row 3
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 622.53

# Row 4
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 618.23

# Row 5
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 592.12

# Row 6
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 470.68

# Row 7
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 468.38

# Row 8
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 381.62

# Row 12
$ws1.Range("B12").Value = 3
$ws1.Range("D12").Value = 18.32

# Row 14
$ws1.Range("A14").Value = "SMB CI (SMBC)"
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 10.86
$ws1.Range("E14").Value = 7.41
$ws1.Range("G14").Value = "➖ Neutre"

# Row 15
$ws1.Range("A15").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B15").Value = 3
$ws1.Range("C15").Value = 1
$ws1.Range("D15").Value = 6.86
$ws1.Range("E15").Value = -4.34
$ws1.Range("F15").Value = "🟢 Achat"
$ws1.Range("G15").Value = "✅ Renforcer"

# Row 16
$ws1.Range("A16").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 6.36
$ws1.Range("E16").Value = 6.36
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

# Row 19
$ws1.Range("A19").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 4.34
$ws1.Range("E19").Value = -1.98
$ws1.Range("G19").Value = "👀 À surveiller"

# Row 20
$ws1.Range("A20").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("C20").Value = 0
$ws1.Range("D20").Value = 3.74
$ws1.Range("E20").Value = 3.74
$ws1.Range("G20").Value = "➖ Neutre"

# Row 21
$ws1.Range("A21").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 3.13
$ws1.Range("E21").Value = 3.13
$ws1.Range("G21").Value = "➖ Neutre"

# Row 22
$ws1.Range("A22").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("D22").Value = 2.81
$ws1.Range("E22").Value = 4.39

# Row 23
$ws1.Range("A23").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B23").Value = 1
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = 1.18
$ws1.Range("E23").Value = -1.38

# Row 24
$ws1.Range("A24").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("D24").Value = 0.18
$ws1.Range("E24").Value = -4.17

# Row 25
$ws1.Range("A25").Value = "SAPH CI (SPHC)"
$ws1.Range("B25").Value = 0
$ws1.Range("D25").Value = -1.31
$ws1.Range("E25").Value = -1.31
$ws1.Range("G25").Value = "➖ Neutre"

# Row 26
$ws1.Range("A26").Value = "FILTISAC CI (FTSC)"
$ws1.Range("D26").Value = -1.67
$ws1.Range("E26").Value = -1.67

# Row 27
$ws1.Range("A27").Value = "CIE CI (CIEC)"
$ws1.Range("B27").Value = 0
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = -1.85
$ws1.Range("E27").Value = -1.85
$ws1.Range("G27").Value = "➖ Neutre"

# Row 28
$ws1.Range("A28").Value = "SICABLE CI (CABC)"
$ws1.Range("D28").Value = -1.99
$ws1.Range("E28").Value = -1.99

# Row 29
$ws1.Range("A29").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("D29").Value = -2.08
$ws1.Range("E29").Value = -2.08

# Row 30
$ws1.Range("A30").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("D30").Value = -2.93
$ws1.Range("E30").Value = -2.93

# Row 31
$ws1.Range("A31").Value = "UNIWAX CI (UNXC)"
$ws1.Range("D31").Value = -5.52
$ws1.Range("E31").Value = 5.61

# Update Top_YTD sheet (sheet 2) progression figures
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 6082.71
$ws2.Range("B3").Value = 4168.01
$ws2.Range("B4").Value = 4098.79
$ws2.Range("B5").Value = 3684.26
$ws2.Range("B6").Value = 2144.84
$ws2.Range("B7").Value = 2121.06
$ws2.Range("B8").Value = 1357.82
